$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6
$ws.Range("H6").Value = 288788.12
$ws.Range("I6").Value = 368367.66
$ws.Range("J6").Value = 50049.5
$ws.Range("K6").Value = 1105102.98
$ws.Range("L6").Value = 150148.5
$ws.Range("M6").Value = -1104990.98
$ws.Range("N6").Value = -150372.5
# row 11
$ws.Range("H11").Value = 41985.055
$ws.Range("I11").Value = 41985.055
$ws.Range("K11").Value = 41985.055
$ws.Range("M11").Value = -41845.055
# row 70
$ws.Range("H70").Value = 2584.5
$ws.Range("I70").Value = 2824.625
$ws.Range("K70").Value = 8473.875
$ws.Range("M70").Value = -8203.875
# row 73
$ws.Range("H73").Value = 2584.5
$ws.Range("I73").Value = 2824.625
$ws.Range("K73").Value = 8473.875
$ws.Range("M73").Value = -7537.875
# row 98
$ws.Range("H98").Value = 848.8889
$ws.Range("I98").Value = 848.8889
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 848.8889
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 649.1111
$ws.Range("N98").ClearContents()
# row 112
$ws.Range("H112").Value = 2856.7
$ws.Range("I112").Value = 2188.8333
$ws.Range("K112").Value = 6566.499899999999
$ws.Range("M112").Value = -5458.499899999999
# row 113
$ws.Range("H113").Value = 3977
$ws.Range("I113").Value = 3945
$ws.Range("J113").Value = 3998.3333
$ws.Range("K113").Value = 3945
$ws.Range("L113").Value = 3998.3333
$ws.Range("M113").Value = -691
$ws.Range("N113").Value = -10506.3333
# row 122
$ws.Range("H122").Value = 848.8889
$ws.Range("I122").Value = 848.8889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2546.6667
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -96.66670000000022
$ws.Range("N122").ClearContents()
# row 138
$ws.Range("H138").Value = 4617.44
$ws.Range("I138").Value = 2133.3
$ws.Range("J138").Value = 6273.533
$ws.Range("K138").Value = 6399.900000000001
$ws.Range("L138").Value = 18820.599
$ws.Range("M138").Value = -1259.900000000001
$ws.Range("N138").Value = -29100.599
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 660.7568
$ws.Range("J2").Value = 997.875
$ws.Range("L2").Value = 997.875
$ws.Range("N2").Value = -1223.875
# row 74
$ws.Range("H74").Value = 2534700.5
$ws.Range("I74").Value = 3707461
$ws.Range("K74").Value = 3707461
$ws.Range("M74").Value = -3706587
# row 77
$ws.Range("H77").Value = 2534700.5
$ws.Range("I77").Value = 3707461
$ws.Range("K77").Value = 18537305
$ws.Range("M77").Value = -18532937
# row 116
$ws.Range("H116").Value = 660.7568
$ws.Range("J116").Value = 997.875
$ws.Range("L116").Value = 997.875
$ws.Range("N116").Value = -5585.875
# row 139
$ws.Range("H139").Value = 119373.5
$ws.Range("J139").Value = 119373.5
$ws.Range("L139").Value = 119373.5
$ws.Range("N139").Value = -129653.5
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 660.7568
$ws.Range("J3").Value = 997.875
$ws.Range("L3").Value = 997.875
$ws.Range("N3").Value = -1225.875
# row 94
$ws.Range("H94").Value = 4716.9
$ws.Range("I94").Value = 5221.25
$ws.Range("K94").Value = 5221.25
$ws.Range("M94").Value = -4770.25
# row 138
$ws.Range("H138").Value = 86222.664
$ws.Range("J138").Value = 86222.664
$ws.Range("L138").Value = 86222.664
$ws.Range("N138").Value = -96502.664
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 4634262
$ws.Range("I31").Value = 6948768.5
$ws.Range("J31").Value = 5250
$ws.Range("K31").Value = 6948768.5
$ws.Range("L31").Value = 5250
$ws.Range("M31").Value = -6948473.5
$ws.Range("N31").Value = -5840
# row 34
$ws.Range("H34").Value = 4634262
$ws.Range("I34").Value = 6948768.5
$ws.Range("J34").Value = 5250
$ws.Range("K34").Value = 6948768.5
$ws.Range("L34").Value = 5250
$ws.Range("M34").Value = -6948566.5
$ws.Range("N34").Value = -5654
# row 132
$ws.Range("H132").Value = 2220.577
$ws.Range("I132").Value = 2335.389
$ws.Range("K132").Value = 7006.167
$ws.Range("M132").Value = -4476.167
$ws = $wb.Worksheets.Item("CUL")
# row 7
$ws.Range("H7").Value = 218.33333
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# row 92
$ws.Range("H92").Value = 461.75
$ws.Range("I92").Value = 444.5
$ws.Range("J92").Value = 467.5
$ws.Range("K92").Value = 1333.5
$ws.Range("L92").Value = 1402.5
$ws.Range("M92").Value = -85.5
$ws.Range("N92").Value = -3898.5
# row 107
$ws.Range("H107").Value = 301.82144
$ws.Range("I107").Value = 336.375
$ws.Range("J107").Value = 288
$ws.Range("K107").Value = 1009.125
$ws.Range("L107").Value = 864
$ws.Range("M107").Value = 910.875
$ws.Range("N107").Value = -4704
# row 109
$ws.Range("H109").Value = 2426.3333
$ws.Range("J109").Value = 3625
$ws.Range("L109").Value = 10875
$ws.Range("N109").Value = -12955
# row 117
$ws.Range("H117").Value = 803.36365
$ws.Range("I117").Value = 895.2
$ws.Range("J117").Value = 726.8333
$ws.Range("K117").Value = 2685.6
$ws.Range("L117").Value = 2180.4999
$ws.Range("M117").Value = 756.3999999999996
$ws.Range("N117").Value = -9064.499899999999
# row 131
$ws.Range("H131").Value = 5849894
$ws.Range("I131").Value = 1388.3334
$ws.Range("K131").Value = 4165.0002
$ws.Range("M131").Value = 874.9997999999996
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 2885773.2
$ws.Range("I2").Value = 3740783.8
$ws.Range("K2").Value = 3740783.8
$ws.Range("M2").Value = -3740670.8
# row 113
$ws.Range("H113").Value = 4390.778
$ws.Range("I113").Value = 3772.7693
$ws.Range("J113").Value = 5997.6
$ws.Range("K113").Value = 3772.7693
$ws.Range("L113").Value = 5997.6
$ws.Range("M113").Value = -1602.7693
$ws.Range("N113").Value = -10337.6
# row 122
$ws.Range("H122").Value = 5725.4346
$ws.Range("I122").Value = 6433
$ws.Range("J122").Value = 3720.6667
$ws.Range("K122").Value = 19299
$ws.Range("L122").Value = 11162.0001
$ws.Range("M122").Value = -16849
$ws.Range("N122").Value = -16062.0001
# row 132
$ws.Range("H132").Value = 14252.763
$ws.Range("I132").Value = 13690.533
$ws.Range("K132").Value = 41071.599
$ws.Range("M132").Value = -38541.599
# row 136
$ws.Range("H136").Value = 61428.805
$ws.Range("J136").Value = 61428.805
$ws.Range("L136").Value = 184286.415
$ws.Range("N136").Value = -189386.415
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 3381.6
$ws.Range("J40").Value = 6165
$ws.Range("L40").Value = 6165
$ws.Range("N40").Value = -6437
# row 61
$ws.Range("H61").Value = 16995.889
$ws.Range("I61").Value = 10709.429
$ws.Range("K61").Value = 10709.429
$ws.Range("M61").Value = -10507.429
# row 113
$ws.Range("H113").Value = 16995.889
$ws.Range("I113").Value = 10709.429
$ws.Range("K113").Value = 10709.429
$ws.Range("M113").Value = -8539.429
# row 122
$ws.Range("H122").Value = 3726.4
$ws.Range("J122").Value = 4845.7144
$ws.Range("L122").Value = 14537.1432
$ws.Range("N122").Value = -19437.1432
# row 127
$ws.Range("H127").Value = 43995.75
$ws.Range("J127").Value = 43995.75
$ws.Range("L127").Value = 43995.75
$ws.Range("N127").Value = -53915.75
# row 132
$ws.Range("H132").Value = 4873340.5
$ws.Range("I132").Value = 8348830.5
$ws.Range("J132").Value = 7654.9
$ws.Range("K132").Value = 25046491.5
$ws.Range("L132").Value = 22964.7
$ws.Range("M132").Value = -25043961.5
$ws.Range("N132").Value = -28024.7
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 2778635
$ws.Range("I113").Value = 866.5925999999999
$ws.Range("J113").Value = 8547846
$ws.Range("K113").Value = 2599.7778
$ws.Range("L113").Value = 25643538
$ws.Range("M113").Value = -429.7777999999998
$ws.Range("N113").Value = -25647878
# row 122
$ws.Range("H122").Value = 52332.348
$ws.Range("I122").Value = 3682.2856
$ws.Range("K122").Value = 11046.8568
$ws.Range("M122").Value = -8596.856800000001
